# Modificada lectura de fichero maestro y fichero de datos
#
# The "Tipo" (Type) column (E) used to hold bare numeric codes (1/2/3).
# It is changed here to hold the actual text labels the codes stood for:
#   1 -> Person, 2 -> Entity, 3 -> Sensor
# so the loader can read the type directly instead of re-mapping a code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Person"
$ws.Range("E3").Value = "Person"
$ws.Range("E4").Value = "Person"
$ws.Range("E5").Value = "Sensor"
$ws.Range("E6").Value = "Entity"
$ws.Range("E7").Value = "Sensor"

# Row 4's location cell picked up an underline (black, not the hyperlink
# style) along the way.
$ws.Range("B4").Font.Underline = $true

# Leave the selection where the author left it when they saved.
$ws.Range("E7").Select()
